$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "5e jour" (column F) time entries added for newly-added tasks
# (register + login-by-mail related rows) in the "Planning effectif" sheet.
$ws.Range("F9").Value = 0.020833333333333332
$ws.Range("F10").Value = 0.020833333333333332
$ws.Range("F11").Value = 0.020833333333333332
$ws.Range("F19").Value = 0.10416666666666667
$ws.Range("F22").Value = 0.020833333333333332
$ws.Range("F23").Value = 0.020833333333333332
$ws.Range("F24").Value = 0.020833333333333332
$ws.Range("F38").Value = 0.10416666666666667
